# WMT.xlsx update: "Cost of Debt" assumption row removed from the debt
# paydown schedule on Sheet2; the interest-rate driver used by the
# I9:R9 "Interest Net" projection now points at U16 (the ROIC/"Discount"
# style row) whose value is lowered from 4% to 2%, instead of the old,
# now-deleted U15 "Cost of Debt" assumption (which also held 2%).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Remove the "Cost of Debt" label (T15) entirely, and blank out its
# value cell (U15) while keeping its number-format style.
$ws.Range("T15").Clear()
$ws.Range("U15").ClearContents()

# The remaining rate row (T16 "ROIC"/U16) drops from 4% to 2%.
$ws.Range("U16").Value = 0.02

# Re-point the Interest Net projection formulas from the deleted
# $U$15 onto $U$16.
$ws.Range("I9").Formula = "=H26*`$U`$16"
$ws.Range("J9:R9").Formula = "=I26*`$U`$16"

# Restore the user's on-screen selection.
$ws.Range("S21").Select()
